# Applies updated metric values for experiments 20, 22, 23 on both
# "Sheet1" (detailed per-experiment metrics, columns M:AB) and
# "Sheet2" (summary metrics, columns G:L).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1 - row 20
# ---------------------------------------------------------------------
$ws1.Range("M20").Value  = 0.03
$ws1.Range("N20").Value  = 0.114
$ws1.Range("O20").Value  = 0.078
$ws1.Range("P20").Value  = 0.067
$ws1.Range("Q20").Value  = 0.99
$ws1.Range("R20").Value  = 0.972
$ws1.Range("S20").Value  = 0.981
$ws1.Range("T20").Value  = 0.981
$ws1.Range("U20").Value  = 0.026
$ws1.Range("V20").Value  = 0.091
$ws1.Range("W20").Value  = 0.078
$ws1.Range("X20").Value  = 0.063
$ws1.Range("Y20").Value  = 0.991
$ws1.Range("Z20").Value  = 0.976
$ws1.Range("AA20").Value = 0.981
$ws1.Range("AB20").Value = 0.982

# ---------------------------------------------------------------------
# Sheet1 - row 22
# ---------------------------------------------------------------------
$ws1.Range("M22").Value  = 0.034
$ws1.Range("N22").Value  = 0.108
$ws1.Range("O22").Value  = 0.075
$ws1.Range("P22").Value  = 0.064
$ws1.Range("Q22").Value  = 0.989
$ws1.Range("R22").Value  = 0.972
$ws1.Range("S22").Value  = 0.982
$ws1.Range("T22").Value  = 0.98
$ws1.Range("U22").Value  = 0.024
$ws1.Range("V22").Value  = 0.091
$ws1.Range("W22").Value  = 0.075
$ws1.Range("X22").Value  = 0.06
$ws1.Range("Y22").Value  = 0.992
$ws1.Range("Z22").Value  = 0.975
$ws1.Range("AA22").Value = 0.982
$ws1.Range("AB22").Value = 0.982

# ---------------------------------------------------------------------
# Sheet1 - row 23
# ---------------------------------------------------------------------
$ws1.Range("N23").Value  = 0.093
$ws1.Range("O23").Value  = 0.093
$ws1.Range("P23").Value  = 0.093
$ws1.Range("S23").Value  = 0.981
$ws1.Range("T23").Value  = 0.98
$ws1.Range("W23").Value  = 0.091
$ws1.Range("X23").Value  = 0.08799999999999999
$ws1.Range("AA23").Value = 0.981
$ws1.Range("AB23").Value = 0.981

# ---------------------------------------------------------------------
# Sheet2 - row 20
# ---------------------------------------------------------------------
$ws2.Range("G20").Value = 0.972
$ws2.Range("H20").Value = 0.976
$ws2.Range("I20").Value = 0.981
$ws2.Range("J20").Value = 0.981
$ws2.Range("K20").Value = 0.981
$ws2.Range("L20").Value = 0.982

# ---------------------------------------------------------------------
# Sheet2 - row 22
# ---------------------------------------------------------------------
$ws2.Range("G22").Value = 0.972
$ws2.Range("H22").Value = 0.975
$ws2.Range("I22").Value = 0.982
$ws2.Range("J22").Value = 0.982
$ws2.Range("K22").Value = 0.98
$ws2.Range("L22").Value = 0.982

# ---------------------------------------------------------------------
# Sheet2 - row 23
# ---------------------------------------------------------------------
$ws2.Range("I23").Value = 0.981
$ws2.Range("J23").Value = 0.981
$ws2.Range("K23").Value = 0.98
$ws2.Range("L23").Value = 0.981
